$wb = $excel.ActiveWorkbook

# --- paper sheet (sheetId 13) ---
$wsPaper = $wb.Worksheets.Item("paper")

$wsPaper.Range("E2").Value = -261.092183753968
$wsPaper.Range("F2").Value = 0.134484074492
$wsPaper.Range("G2").Value = -23431.99349141002
$wsPaper.Range("H2").Value = 23.258877643235
$wsPaper.Range("I2").Value = -0.005769391894
$wsPaper.Range("J2").Value = -261.092183753968
$wsPaper.Range("K2").Value = 0.134484074492
$wsPaper.Range("L2").Value = -23431.99349141002
$wsPaper.Range("M2").Value = 23.258877643235
$wsPaper.Range("N2").Value = -0.005769391894
$wsPaper.Range("E3").Value = -30.44046842074
$wsPaper.Range("F3").Value = 0.025926991566
$wsPaper.Range("G3").Value = 73472.61084267939
$wsPaper.Range("H3").Value = -73.32959183150901
$wsPaper.Range("I3").Value = 0.01830174419
$wsPaper.Range("J3").Value = -30.44046842074
$wsPaper.Range("K3").Value = 0.025926991566
$wsPaper.Range("L3").Value = 73472.61084267939
$wsPaper.Range("M3").Value = -73.32959183150901
$wsPaper.Range("N3").Value = 0.01830174419
$wsPaper.Range("E4").Value = -849.700794374014
$wsPaper.Range("F4").Value = 0.433384170937
$wsPaper.Range("G4").Value = -6148.130282926714
$wsPaper.Range("H4").Value = 5.721178856293
$wsPaper.Range("I4").Value = -0.001319271777
$wsPaper.Range("J4").Value = -849.700794374014
$wsPaper.Range("K4").Value = 0.433384170937
$wsPaper.Range("L4").Value = -6148.130282926714
$wsPaper.Range("M4").Value = 5.721178856293
$wsPaper.Range("N4").Value = -0.001319271777
$wsPaper.Range("E5").Value = 587.071660179378
$wsPaper.Range("F5").Value = -0.289052557496
$wsPaper.Range("G5").Value = 6241.232902265833
$wsPaper.Range("H5").Value = -5.931864985338
$wsPaper.Range("I5").Value = 0.001407846488
$wsPaper.Range("J5").Value = 587.071660179378
$wsPaper.Range("K5").Value = -0.289052557496
$wsPaper.Range("L5").Value = 6241.232902265833
$wsPaper.Range("M5").Value = -5.931864985338
$wsPaper.Range("N5").Value = 0.001407846488
$wsPaper.Range("E6").Value = 5.143238135853
$wsPaper.Range("F6").Value = 0.0007394641779999999
$wsPaper.Range("G6").Value = 10423.79257958418
$wsPaper.Range("H6").Value = -10.396997966379
$wsPaper.Range("I6").Value = 0.002594170604
$wsPaper.Range("J6").Value = 5.143238135853
$wsPaper.Range("K6").Value = 0.0007394641779999999
$wsPaper.Range("L6").Value = 10423.79257958418
$wsPaper.Range("M6").Value = -10.396997966379
$wsPaper.Range("N6").Value = 0.002594170604
$wsPaper.Range("E7").Value = 1429.548678857135
$wsPaper.Range("F7").Value = -0.7061391125649999
$wsPaper.Range("G7").Value = -338458.8935053673
$wsPaper.Range("H7").Value = 338.452596120029
$wsPaper.Range("I7").Value = -0.084606215098
$wsPaper.Range("J7").Value = 1429.548678857135
$wsPaper.Range("K7").Value = -0.7061391125649999
$wsPaper.Range("L7").Value = -338458.8935053673
$wsPaper.Range("M7").Value = 338.452596120029
$wsPaper.Range("N7").Value = -0.084606215098
$wsPaper.Range("E8").Value = 324.334709019177
$wsPaper.Range("F8").Value = -0.158911850001
$wsPaper.Range("G8").Value = -57003.49427617952
$wsPaper.Range("H8").Value = 57.053850929205
$wsPaper.Range("I8").Value = -0.014274227292
$wsPaper.Range("J8").Value = 324.334709019177
$wsPaper.Range("K8").Value = -0.158911850001
$wsPaper.Range("L8").Value = -57003.49427617952
$wsPaper.Range("M8").Value = 57.053850929205
$wsPaper.Range("N8").Value = -0.014274227292
$wsPaper.Range("E9").Value = 1161.572768553227
$wsPaper.Range("F9").Value = -0.572861816244
$wsPaper.Range("G9").Value = 11937.86450335169
$wsPaper.Range("H9").Value = -11.327523794142
$wsPaper.Range("I9").Value = 0.002683221052
$wsPaper.Range("J9").Value = 1161.572768553227
$wsPaper.Range("K9").Value = -0.572861816244
$wsPaper.Range("L9").Value = 11937.86450335169
$wsPaper.Range("M9").Value = -11.327523794142
$wsPaper.Range("N9").Value = 0.002683221052
$wsPaper.Range("E10").Value = 374.276456625246
$wsPaper.Range("F10").Value = -0.182182243025
$wsPaper.Range("G10").Value = 19271.81139776916
$wsPaper.Range("H10").Value = -19.041786781313
$wsPaper.Range("I10").Value = 0.004705353644
$wsPaper.Range("J10").Value = 374.276456625246
$wsPaper.Range("K10").Value = -0.182182243025
$wsPaper.Range("L10").Value = 19271.81139776916
$wsPaper.Range("M10").Value = -19.041786781313
$wsPaper.Range("N10").Value = 0.004705353644
$wsPaper.Range("E11").Value = 777.506289682108
$wsPaper.Range("F11").Value = -0.384286240911
$wsPaper.Range("G11").Value = 125768.3547636789
$wsPaper.Range("H11").Value = -125.124257907391
$wsPaper.Range("I11").Value = 0.031121844523
$wsPaper.Range("J11").Value = 777.506289682108
$wsPaper.Range("K11").Value = -0.384286240911
$wsPaper.Range("L11").Value = 125768.3547636789
$wsPaper.Range("M11").Value = -125.124257907391
$wsPaper.Range("N11").Value = 0.031121844523
$wsPaper.Range("E12").Value = 343.814097762301
$wsPaper.Range("F12").Value = -0.167481491958
$wsPaper.Range("G12").Value = 23454.69964288389
$wsPaper.Range("H12").Value = -23.231979759489
$wsPaper.Range("I12").Value = 0.00575444839
$wsPaper.Range("J12").Value = 343.814097762301
$wsPaper.Range("K12").Value = -0.167481491958
$wsPaper.Range("L12").Value = 23454.69964288389
$wsPaper.Range("M12").Value = -23.231979759489
$wsPaper.Range("N12").Value = 0.00575444839
$wsPaper.Range("E13").Value = 1742.488129895361
$wsPaper.Range("F13").Value = -0.863895917567
$wsPaper.Range("G13").Value = 112264.9637034967
$wsPaper.Range("H13").Value = -111.164535043386
$wsPaper.Range("I13").Value = 0.027519321159
$wsPaper.Range("J13").Value = 1742.488129895361
$wsPaper.Range("K13").Value = -0.863895917567
$wsPaper.Range("L13").Value = 112264.9637034967
$wsPaper.Range("M13").Value = -111.164535043386
$wsPaper.Range("N13").Value = 0.027519321159
$wsPaper.Range("E14").Value = -870.123772240489
$wsPaper.Range("F14").Value = 0.439246136851
$wsPaper.Range("G14").Value = -95803.6069123714
$wsPaper.Range("H14").Value = 95.086536885925
$wsPaper.Range("I14").Value = -0.023590343945
$wsPaper.Range("J14").Value = -870.123772240489
$wsPaper.Range("K14").Value = 0.439246136851
$wsPaper.Range("L14").Value = -95803.6069123714
$wsPaper.Range("M14").Value = 95.086536885925
$wsPaper.Range("N14").Value = -0.023590343945
$wsPaper.Range("E15").Value = 24.043195028042
$wsPaper.Range("F15").Value = -0.007561929633
$wsPaper.Range("G15").Value = -78203.01267452863
$wsPaper.Range("H15").Value = 78.068235462968
$wsPaper.Range("I15").Value = -0.019480869813
$wsPaper.Range("J15").Value = 24.043195028042
$wsPaper.Range("K15").Value = -0.007561929633
$wsPaper.Range("L15").Value = -78203.01267452863
$wsPaper.Range("M15").Value = 78.068235462968
$wsPaper.Range("N15").Value = -0.019480869813
$wsPaper.Range("E16").Value = 254.49426467161
$wsPaper.Range("F16").Value = -0.123276437148
$wsPaper.Range("G16").Value = -29340.02809988912
$wsPaper.Range("H16").Value = 29.411844944209
$wsPaper.Range("I16").Value = -0.007368828478
$wsPaper.Range("J16").Value = 254.49426467161
$wsPaper.Range("K16").Value = -0.123276437148
$wsPaper.Range("L16").Value = -29340.02809988912
$wsPaper.Range("M16").Value = 29.411844944209
$wsPaper.Range("N16").Value = -0.007368828478
$wsPaper.Range("E17").Value = 11.073362943876
$wsPaper.Range("F17").Value = -0.000430313369
$wsPaper.Range("G17").Value = 32622.41541345948
$wsPaper.Range("H17").Value = -32.54631613663
$wsPaper.Range("I17").Value = 0.008119995418999999
$wsPaper.Range("J17").Value = 11.073362943876
$wsPaper.Range("K17").Value = -0.000430313369
$wsPaper.Range("L17").Value = 32622.41541345948
$wsPaper.Range("M17").Value = -32.54631613663
$wsPaper.Range("N17").Value = 0.008119995418999999
$wsPaper.Range("E18").Value = 1625.981122865376
$wsPaper.Range("F18").Value = -0.80247621196
$wsPaper.Range("G18").Value = 142049.7456219967
$wsPaper.Range("H18").Value = -140.944387550103
$wsPaper.Range("I18").Value = 0.034964532359
$wsPaper.Range("J18").Value = 1625.981122865376
$wsPaper.Range("K18").Value = -0.80247621196
$wsPaper.Range("L18").Value = 142049.7456219967
$wsPaper.Range("M18").Value = -140.944387550103
$wsPaper.Range("N18").Value = 0.034964532359
$wsPaper.Range("E19").Value = 1818.988797160781
$wsPaper.Range("F19").Value = -0.892545496068
$wsPaper.Range("G19").Value = 84408.01367514276
$wsPaper.Range("H19").Value = -83.31580087051999
$wsPaper.Range("I19").Value = 0.020564087875
$wsPaper.Range("J19").Value = 1818.988797160781
$wsPaper.Range("K19").Value = -0.892545496068
$wsPaper.Range("L19").Value = 84408.01367514276
$wsPaper.Range("M19").Value = -83.31580087051999
$wsPaper.Range("N19").Value = 0.020564087875
$wsPaper.Range("E20").Value = 1314.543556275201
$wsPaper.Range("F20").Value = -0.646916305793
$wsPaper.Range("G20").Value = -306860.4957595402
$wsPaper.Range("H20").Value = 306.909565963935
$wsPaper.Range("I20").Value = -0.07673342309800001
$wsPaper.Range("J20").Value = 1314.543556275201
$wsPaper.Range("K20").Value = -0.646916305793
$wsPaper.Range("L20").Value = -306860.4957595402
$wsPaper.Range("M20").Value = 306.909565963935
$wsPaper.Range("N20").Value = -0.07673342309800001
$wsPaper.Range("E21").Value = -59.002884210532
$wsPaper.Range("F21").Value = 0.032928102865
$wsPaper.Range("G21").Value = -113808.7627995997
$wsPaper.Range("H21").Value = 113.554373889886
$wsPaper.Range("I21").Value = -0.028322892322
$wsPaper.Range("J21").Value = -59.002884210532
$wsPaper.Range("K21").Value = 0.032928102865
$wsPaper.Range("L21").Value = -113808.7627995997
$wsPaper.Range("M21").Value = 113.554373889886
$wsPaper.Range("N21").Value = -0.028322892322
$wsPaper.Range("E22").Value = -70.279892033198
$wsPaper.Range("F22").Value = 0.044676056489
$wsPaper.Range("G22").Value = 123646.0897581631
$wsPaper.Range("H22").Value = -123.408307251605
$wsPaper.Range("I22").Value = 0.030797011516
$wsPaper.Range("J22").Value = -70.279892033198
$wsPaper.Range("K22").Value = 0.044676056489
$wsPaper.Range("L22").Value = 123646.0897581631
$wsPaper.Range("M22").Value = -123.408307251605
$wsPaper.Range("N22").Value = 0.030797011516
$wsPaper.Range("E23").Value = -681.675656407076
$wsPaper.Range("F23").Value = 0.346602361737
$wsPaper.Range("G23").Value = 16076.27777509058
$wsPaper.Range("H23").Value = -16.377715152342
$wsPaper.Range("I23").Value = 0.004172612856
$wsPaper.Range("J23").Value = -681.675656407076
$wsPaper.Range("K23").Value = 0.346602361737
$wsPaper.Range("L23").Value = 16076.27777509058
$wsPaper.Range("M23").Value = -16.377715152342
$wsPaper.Range("N23").Value = 0.004172612856
$wsPaper.Range("E24").Value = 165.06117645879
$wsPaper.Range("F24").Value = -0.074595799711
$wsPaper.Range("G24").Value = 55908.97956746679
$wsPaper.Range("H24").Value = -55.706627148947
$wsPaper.Range("I24").Value = 0.013879844664
$wsPaper.Range("J24").Value = 165.06117645879
$wsPaper.Range("K24").Value = -0.074595799711
$wsPaper.Range("L24").Value = 55908.97956746679
$wsPaper.Range("M24").Value = -55.706627148947
$wsPaper.Range("N24").Value = 0.013879844664
$wsPaper.Range("E25").Value = -194.887075685634
$wsPaper.Range("F25").Value = 0.102470383992
$wsPaper.Range("G25").Value = 69427.78368029829
$wsPaper.Range("H25").Value = -69.380456434069
$wsPaper.Range("I25").Value = 0.017335556651
$wsPaper.Range("J25").Value = -194.887075685634
$wsPaper.Range("K25").Value = 0.102470383992
$wsPaper.Range("L25").Value = 69427.78368029829
$wsPaper.Range("M25").Value = -69.380456434069
$wsPaper.Range("N25").Value = 0.017335556651
$wsPaper.Range("E26").Value = 236.474722555805
$wsPaper.Range("F26").Value = -0.114006994389
$wsPaper.Range("G26").Value = -76893.97542606578
$wsPaper.Range("H26").Value = 76.861629891246
$wsPaper.Range("I26").Value = -0.019204941056
$wsPaper.Range("J26").Value = 236.474722555805
$wsPaper.Range("K26").Value = -0.114006994389
$wsPaper.Range("L26").Value = -76893.97542606578
$wsPaper.Range("M26").Value = 76.861629891246
$wsPaper.Range("N26").Value = -0.019204941056
$wsPaper.Range("E27").Value = 7.470796460177
$wsPaper.Range("J27").Value = 7.470796460177
$wsPaper.Range("E28").Value = 7.470796460177
$wsPaper.Range("J28").Value = 7.470796460177
$wsPaper.Range("E29").Value = 318.094471451535
$wsPaper.Range("F29").Value = -0.156444422204
$wsPaper.Range("G29").Value = -154410.2332071213
$wsPaper.Range("H29").Value = 154.242473911737
$wsPaper.Range("I29").Value = -0.03851689233
$wsPaper.Range("J29").Value = 318.094471451535
$wsPaper.Range("K29").Value = -0.156444422204
$wsPaper.Range("L29").Value = -154410.2332071213
$wsPaper.Range("M29").Value = 154.242473911737
$wsPaper.Range("N29").Value = -0.03851689233
$wsPaper.Range("E30").Value = 7.470796460177
$wsPaper.Range("J30").Value = 7.470796460177
$wsPaper.Range("E31").Value = -741.764204484433
$wsPaper.Range("F31").Value = 0.371006100729
$wsPaper.Range("G31").Value = -143314.1684743427
$wsPaper.Range("H31").Value = 142.349311737595
$wsPaper.Range("I31").Value = -0.035346310777
$wsPaper.Range("J31").Value = -741.764204484433
$wsPaper.Range("K31").Value = 0.371006100729
$wsPaper.Range("L31").Value = -143314.1684743427
$wsPaper.Range("M31").Value = 142.349311737595
$wsPaper.Range("N31").Value = -0.035346310777
$wsPaper.Range("E32").Value = 7.470796460177
$wsPaper.Range("J32").Value = 7.470796460177
$wsPaper.Range("E34").Value = 4114.612607193242
$wsPaper.Range("F34").Value = -2.042350397733
$wsPaper.Range("G34").Value = 512393.0377879965
$wsPaper.Range("H34").Value = -509.30057873856
$wsPaper.Range("I34").Value = 0.126557762555
$wsPaper.Range("J34").Value = 4114.612607193242
$wsPaper.Range("K34").Value = -2.042350397733
$wsPaper.Range("L34").Value = 512393.0377879965
$wsPaper.Range("M34").Value = -509.30057873856
$wsPaper.Range("N34").Value = 0.126557762555
$wsPaper.Range("E35").Value = 1174.185759398031
$wsPaper.Range("F35").Value = -0.5753745364940001
$wsPaper.Range("G35").Value = 21637.29844283697
$wsPaper.Range("H35").Value = -20.997414449849
$wsPaper.Range("I35").Value = 0.00509517152
$wsPaper.Range("J35").Value = 1174.185759398031
$wsPaper.Range("K35").Value = -0.5753745364940001
$wsPaper.Range("L35").Value = 21637.29844283697
$wsPaper.Range("M35").Value = -20.997414449849
$wsPaper.Range("N35").Value = 0.00509517152

# --- glass sheet (sheetId 15) ---
$wsGlass = $wb.Worksheets.Item("glass")

$wsGlass.Range("E2").Value = 6.133162642233
$wsGlass.Range("J2").Value = 6.133162642233
$wsGlass.Range("E3").Value = 6.133162642233
$wsGlass.Range("J3").Value = 6.133162642233
$wsGlass.Range("E4").Value = 6.133162642233
$wsGlass.Range("J4").Value = 6.133162642233
$wsGlass.Range("E5").Value = 6.133162642233
$wsGlass.Range("J5").Value = 6.133162642233
$wsGlass.Range("E6").Value = 6.133162642233
$wsGlass.Range("J6").Value = 6.133162642233
$wsGlass.Range("E7").Value = 6.133162642233
$wsGlass.Range("J7").Value = 6.133162642233
$wsGlass.Range("E8").Value = 6.133162642233
$wsGlass.Range("J8").Value = 6.133162642233
$wsGlass.Range("E9").Value = 6.133162642233
$wsGlass.Range("J9").Value = 6.133162642233
$wsGlass.Range("E10").Value = 6.133162642233
$wsGlass.Range("J10").Value = 6.133162642233
$wsGlass.Range("E11").Value = 6.133162642233
$wsGlass.Range("J11").Value = 6.133162642233
$wsGlass.Range("E12").Value = 6.133162642233
$wsGlass.Range("J12").Value = 6.133162642233
$wsGlass.Range("E13").Value = 6.133162642233
$wsGlass.Range("J13").Value = 6.133162642233
$wsGlass.Range("E14").Value = 6.133162642233
$wsGlass.Range("J14").Value = 6.133162642233
$wsGlass.Range("E15").Value = 6.133162642233
$wsGlass.Range("J15").Value = 6.133162642233
$wsGlass.Range("E16").Value = 6.133162642233
$wsGlass.Range("J16").Value = 6.133162642233
$wsGlass.Range("E17").Value = 6.133162642233
$wsGlass.Range("J17").Value = 6.133162642233
$wsGlass.Range("E18").Value = 6.133162642233
$wsGlass.Range("J18").Value = 6.133162642233
$wsGlass.Range("E19").Value = 6.133162642233
$wsGlass.Range("J19").Value = 6.133162642233
$wsGlass.Range("E20").Value = 6.133162642233
$wsGlass.Range("J20").Value = 6.133162642233
$wsGlass.Range("E21").Value = 6.133162642233
$wsGlass.Range("J21").Value = 6.133162642233
$wsGlass.Range("E22").Value = 6.133162642233
$wsGlass.Range("J22").Value = 6.133162642233
$wsGlass.Range("E23").Value = 6.133162642233
$wsGlass.Range("J23").Value = 6.133162642233
$wsGlass.Range("E24").Value = 6.133162642233
$wsGlass.Range("J24").Value = 6.133162642233
$wsGlass.Range("E25").Value = 6.133162642233
$wsGlass.Range("J25").Value = 6.133162642233
$wsGlass.Range("E26").Value = 6.133162642233
$wsGlass.Range("J26").Value = 6.133162642233
$wsGlass.Range("E27").Value = 6.133162642233
$wsGlass.Range("J27").Value = 6.133162642233
$wsGlass.Range("E28").Value = 6.133162642233
$wsGlass.Range("J28").Value = 6.133162642233
$wsGlass.Range("E29").Value = 6.133162642233
$wsGlass.Range("J29").Value = 6.133162642233
$wsGlass.Range("E30").Value = 6.133162642233
$wsGlass.Range("J30").Value = 6.133162642233
$wsGlass.Range("E31").Value = 6.133162642233
$wsGlass.Range("J31").Value = 6.133162642233
$wsGlass.Range("E32").Value = 6.133162642233
$wsGlass.Range("J32").Value = 6.133162642233
$wsGlass.Range("E33").Value = 6.133162642233
$wsGlass.Range("J33").Value = 6.133162642233
$wsGlass.Range("E34").Value = 6.133162642233
$wsGlass.Range("J34").Value = 6.133162642233
$wsGlass.Range("E35").Value = 6.133162642233
$wsGlass.Range("J35").Value = 6.133162642233
